$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: new header row (idx, idx2, Name, Date Start, Date End, (m3/s),
# (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year)
# Columns A-E keep the default (unstyled) look; F-K get the 9pt Arial font
# that the rest of the data table uses.
# ---------------------------------------------------------------------------
# A1:E1 must end up with the plain default style - clear away any leftover
# formatting from the cells that used to live here (e.g. old E1 header).
$ws.Range("A1:E1").ClearFormats()

$ws.Cells.Item(1, 1).Value = "idx"
$ws.Cells.Item(1, 2).Value = "idx2"
$ws.Cells.Item(1, 3).Value = "Name"
$ws.Cells.Item(1, 4).Value = "Date Start"
$ws.Cells.Item(1, 5).Value = "Date End"

$ws.Cells.Item(1, 6).Value = "(m3/s)"
$ws.Cells.Item(1, 6).Font.Size = 9
$ws.Cells.Item(1, 6).NumberFormat = "General"

$ws.Cells.Item(1, 7).Value = "(MW1)"
$ws.Cells.Item(1, 7).Font.Size = 9
$ws.Cells.Item(1, 7).NumberFormat = "General"

$ws.Cells.Item(1, 8).Value = "(MW2)"
$ws.Cells.Item(1, 8).Font.Size = 9
$ws.Cells.Item(1, 8).NumberFormat = "General"

$ws.Cells.Item(1, 9).Value = "(GWh) Winter"
$ws.Cells.Item(1, 9).Font.Size = 9
$ws.Cells.Item(1, 9).NumberFormat = "General"

$ws.Cells.Item(1, 10).Value = "(GWh) Summer"
$ws.Cells.Item(1, 10).Font.Size = 9
$ws.Cells.Item(1, 10).NumberFormat = "General"

$ws.Cells.Item(1, 11).Value = "(GWh) Year"
$ws.Cells.Item(1, 11).Font.Size = 9
$ws.Cells.Item(1, 11).NumberFormat = "General"

# ---------------------------------------------------------------------------
# Data rows: the old table (rows 3-11) slides up by one row (rows 2-10),
# since the two old header rows are merged into a single new header row.
# idx / idx2 / Date Start / Date End -> integer format (style "0")
# Name -> text (9pt font)
# (m3/s) / (MW1) / (MW2) / (GWh) * -> 2 decimal number format (style "0.00")
# ---------------------------------------------------------------------------
$data = @(
    @(1, 304100, "Rotzloch",               1872, 1935, 0.82, 0.38,  0.34,  0.53,  0.65,  1.18),
    @(2, 303000, "Obermatt",                1905, 1963, 11,   7.73,  7.19,  6.29,  23.96, 30.25),
    @(3, 303300, "Oberrickenbach",          1937, 1991, 1,    8.7,   6.8,   4.4,   9.8,   14.2),
    @(4, 303400, "Wolfenschiessen",         1945, 1983, 2.6,  6.6,   6.6,   6.6,   13.7,  20.3),
    @(5, 302600, "Sustli",                  1957, 1998, 0.51, 1.63,  1.56,  3.5,   5.3,   8.8),
    @(6, 303200, "Dallenwil",               1962, 1987, 14.7, 12.22, 11.06, 13.13, 40.67, 53.8),
    @(7, 303100, "Obermatt-Nebenzentrale",  1963, $null, 11,  0.48,  0.44,  0.24,  0.97,  1.21),
    @(8, 302900, "Arni",                    1966, $null, 1,   1.85,  1.74,  0.29,  4.06,  4.35),
    @(9, 302800, "Engelberg",               1967, $null, 1.4, 7.74,  7.64,  2.91,  14.65, 17.56)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Font.Size = 9
    $ws.Cells.Item($r, 1).NumberFormat = "0"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 2).Font.Size = 9
    $ws.Cells.Item($r, 2).NumberFormat = "0"

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 3).Font.Size = 9

    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).Font.Size = 9
    $ws.Cells.Item($r, 4).NumberFormat = "0"

    if ($row[4] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 5).Font.Size = 9
        $ws.Cells.Item($r, 5).NumberFormat = "0"
    } else {
        $ws.Cells.Item($r, 5).Clear()
    }

    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).Font.Size = 9
    $ws.Cells.Item($r, 6).NumberFormat = "0.00"

    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).Font.Size = 9
    $ws.Cells.Item($r, 7).NumberFormat = "0.00"

    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 8).Font.Size = 9
    $ws.Cells.Item($r, 8).NumberFormat = "0.00"

    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 9).Font.Size = 9
    $ws.Cells.Item($r, 9).NumberFormat = "0.00"

    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 10).Font.Size = 9
    $ws.Cells.Item($r, 10).NumberFormat = "0.00"

    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 11).Font.Size = 9
    $ws.Cells.Item($r, 11).NumberFormat = "0.00"

    $r = $r + 1
}

# Remove the now unused old row 11 (table shrank from 9 to 9 data rows + 1
# header instead of 9 data rows + 2 headers: old range was A1:K11, new is
# A1:K10).
$ws.Rows.Item(11).Clear()

# Update the selection to match the target workbook state.
$ws.Range("A2:K2").Select()
